$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -17.34519731065735
$ws.Range("C2").Value = 1.843335999648392
$ws.Range("D2").Value = -17.34519731065735
$ws.Range("E2").Value = -17.34519731065735
$ws.Range("F2").Value = -17.34519731065735
$ws.Range("G2").Value = -17.34519731065735
$ws.Range("H2").Value = -17.34519731065735
$ws.Range("I2").Value = -17.34519731065735
$ws.Range("J2").Value = -17.34519731065735
$ws.Range("K2").Value = -17.34519731065735
$ws.Range("B3").Value = -17.34519731065735
$ws.Range("C3").Value = -17.34519731065735
$ws.Range("D3").Value = -17.34519731065735
$ws.Range("E3").Value = -17.34519731065735
$ws.Range("F3").Value = -17.34519731065735
$ws.Range("G3").Value = -17.34519731065735
$ws.Range("H3").Value = -17.34519731065735
$ws.Range("I3").Value = 2.869393361983968
$ws.Range("J3").Value = -17.34519731065735
$ws.Range("K3").Value = -17.34519731065735
$ws.Range("B4").Value = -17.34519731065735
$ws.Range("C4").Value = 2.204273402820893
$ws.Range("D4").Value = 2.114084206221886
$ws.Range("E4").Value = -17.34519731065735
$ws.Range("F4").Value = 3.502221771365214
$ws.Range("G4").Value = -17.34519731065735
$ws.Range("H4").Value = 1.629242487549975
$ws.Range("I4").Value = -17.34519731065735
$ws.Range("J4").Value = 2.505979830944474
$ws.Range("K4").Value = -17.34519731065735
$ws.Range("B5").Value = -17.34519731065735
$ws.Range("C5").Value = 1.765859378503801
$ws.Range("D5").Value = -17.34519731065735
$ws.Range("E5").Value = -17.34519731065735
$ws.Range("F5").Value = -17.34519731065735
$ws.Range("G5").Value = 2.818738535298087
$ws.Range("H5").Value = -17.34519731065735
$ws.Range("I5").Value = -17.34519731065735
$ws.Range("J5").Value = -17.34519731065735
$ws.Range("K5").Value = -17.34519731065735
$ws.Range("B6").Value = -17.34519731065735
$ws.Range("C6").Value = -17.34519731065735
$ws.Range("D6").Value = -17.34519731065735
$ws.Range("E6").Value = -17.34519731065735
$ws.Range("F6").Value = -17.34519731065735
$ws.Range("G6").Value = -17.34519731065735
$ws.Range("H6").Value = -17.34519731065735
$ws.Range("I6").Value = -17.34519731065735
$ws.Range("J6").Value = -17.34519731065735
$ws.Range("K6").Value = -17.34519731065735
$ws.Range("B7").Value = 2.600585132030595
$ws.Range("C7").Value = -17.34519731065735
$ws.Range("D7").Value = -17.34519731065735
$ws.Range("E7").Value = -17.34519731065735
$ws.Range("F7").Value = -17.34519731065735
$ws.Range("G7").Value = -17.34519731065735
$ws.Range("H7").Value = -17.34519731065735
$ws.Range("I7").Value = -17.34519731065735
$ws.Range("J7").Value = -17.34519731065735
$ws.Range("K7").Value = -17.34519731065735
$ws.Range("B8").Value = -17.34519731065735
$ws.Range("C8").Value = -17.34519731065735
$ws.Range("D8").Value = -17.34519731065735
$ws.Range("E8").Value = 1.837771635717179
$ws.Range("F8").Value = -17.34519731065735
$ws.Range("G8").Value = -17.34519731065735
$ws.Range("H8").Value = -17.34519731065735
$ws.Range("I8").Value = -17.34519731065735
$ws.Range("J8").Value = -17.34519731065735
$ws.Range("K8").Value = -17.34519731065735
$ws.Range("B9").Value = 3.800596178606936
$ws.Range("C9").Value = -17.34519731065735
$ws.Range("D9").Value = -17.34519731065735
$ws.Range("E9").Value = -17.34519731065735
$ws.Range("F9").Value = -17.34519731065735
$ws.Range("G9").Value = -17.34519731065735
$ws.Range("H9").Value = -17.34519731065735
$ws.Range("I9").Value = -17.34519731065735
$ws.Range("J9").Value = -17.34519731065735
$ws.Range("K9").Value = -17.34519731065735
$ws.Range("B10").Value = -17.34519731065735
$ws.Range("C10").Value = -17.34519731065735
$ws.Range("D10").Value = -17.34519731065735
$ws.Range("E10").Value = -17.34519731065735
$ws.Range("F10").Value = -17.34519731065735
$ws.Range("G10").Value = -17.34519731065735
$ws.Range("H10").Value = -17.34519731065735
$ws.Range("I10").Value = 0.862145949511402
$ws.Range("J10").Value = -17.34519731065735
$ws.Range("K10").Value = -17.34519731065735
$ws.Range("B11").Value = -17.34519731065735
$ws.Range("C11").Value = -17.34519731065735
$ws.Range("D11").Value = -17.34519731065735
$ws.Range("E11").Value = 2.895944452862371
$ws.Range("F11").Value = -17.34519731065735
$ws.Range("G11").Value = 2.799178381375301
$ws.Range("H11").Value = -17.34519731065735
$ws.Range("I11").Value = -17.34519731065735
$ws.Range("J11").Value = -17.34519731065735
$ws.Range("K11").Value = -17.34519731065735
$ws.Range("B12").Value = -17.34519731065735
$ws.Range("C12").Value = -17.34519731065735
$ws.Range("D12").Value = -17.34519731065735
$ws.Range("E12").Value = -17.34519731065735
$ws.Range("F12").Value = -17.34519731065735
$ws.Range("G12").Value = -17.34519731065735
$ws.Range("H12").Value = -17.34519731065735
$ws.Range("I12").Value = -17.34519731065735
$ws.Range("J12").Value = -17.34519731065735
$ws.Range("K12").Value = -17.34519731065735
$ws.Range("B13").Value = -17.34519731065735
$ws.Range("C13").Value = -17.34519731065735
$ws.Range("D13").Value = -17.34519731065735
$ws.Range("E13").Value = 2.384464918809092
$ws.Range("F13").Value = -17.34519731065735
$ws.Range("G13").Value = -17.34519731065735
$ws.Range("H13").Value = -17.34519731065735
$ws.Range("I13").Value = -17.34519731065735
$ws.Range("J13").Value = 1.628906339304195
$ws.Range("K13").Value = 4.321919863468304
$ws.Range("B14").Value = -17.34519731065735
$ws.Range("C14").Value = -17.34519731065735
$ws.Range("D14").Value = 1.196583726533412
$ws.Range("E14").Value = -17.34519731065735
$ws.Range("F14").Value = -17.34519731065735
$ws.Range("G14").Value = -17.34519731065735
$ws.Range("H14").Value = -17.34519731065735
$ws.Range("I14").Value = -17.34519731065735
$ws.Range("J14").Value = -17.34519731065735
$ws.Range("K14").Value = -17.34519731065735
$ws.Range("B15").Value = -17.34519731065735
$ws.Range("C15").Value = -17.34519731065735
$ws.Range("D15").Value = 1.200208391132012
$ws.Range("E15").Value = -17.34519731065735
$ws.Range("F15").Value = -17.34519731065735
$ws.Range("G15").Value = -17.34519731065735
$ws.Range("H15").Value = -17.34519731065735
$ws.Range("I15").Value = -17.34519731065735
$ws.Range("J15").Value = -17.34519731065735
$ws.Range("K15").Value = -17.34519731065735
$ws.Range("B16").Value = -17.34519731065735
$ws.Range("C16").Value = -17.34519731065735
$ws.Range("D16").Value = -17.34519731065735
$ws.Range("E16").Value = -17.34519731065735
$ws.Range("F16").Value = -17.34519731065735
$ws.Range("G16").Value = -17.34519731065735
$ws.Range("H16").Value = -17.34519731065735
$ws.Range("I16").Value = -17.34519731065735
$ws.Range("J16").Value = 2.009897798259155
$ws.Range("K16").Value = -17.34519731065735
$ws.Range("B17").Value = -17.34519731065735
$ws.Range("C17").Value = 2.170373566827952
$ws.Range("D17").Value = 2.491353927473387
$ws.Range("E17").Value = -17.34519731065735
$ws.Range("F17").Value = -17.34519731065735
$ws.Range("G17").Value = -17.34519731065735
$ws.Range("H17").Value = 1.15438779295229
$ws.Range("I17").Value = 1.800000985439891
$ws.Range("J17").Value = 2.110155999816407
$ws.Range("K17").Value = -17.34519731065735
$ws.Range("B18").Value = -17.34519731065735
$ws.Range("C18").Value = -17.34519731065735
$ws.Range("D18").Value = -17.34519731065735
$ws.Range("E18").Value = -17.34519731065735
$ws.Range("F18").Value = -17.34519731065735
$ws.Range("G18").Value = -17.34519731065735
$ws.Range("H18").Value = 1.665584822334097
$ws.Range("I18").Value = 1.273005348016018
$ws.Range("J18").Value = 1.526999597712495
$ws.Range("K18").Value = -17.34519731065735
$ws.Range("B19").Value = -17.34519731065735
$ws.Range("C19").Value = -17.34519731065735
$ws.Range("D19").Value = 1.524532068032685
$ws.Range("E19").Value = -17.34519731065735
$ws.Range("F19").Value = -17.34519731065735
$ws.Range("G19").Value = -17.34519731065735
$ws.Range("H19").Value = 1.681231572933039
$ws.Range("I19").Value = 1.484342722683286
$ws.Range("J19").Value = -17.34519731065735
$ws.Range("K19").Value = -17.34519731065735
$ws.Range("B20").Value = -17.34519731065735
$ws.Range("C20").Value = 0.8978325225840842
$ws.Range("D20").Value = 1.367977730631958
$ws.Range("E20").Value = -17.34519731065735
$ws.Range("F20").Value = 3.11582342769773
$ws.Range("G20").Value = -17.34519731065735
$ws.Range("H20").Value = 1.927106697671963
$ws.Range("I20").Value = 1.122905640208963
$ws.Range("J20").Value = -17.34519731065735
$ws.Range("K20").Value = -17.34519731065735
$ws.Range("B21").Value = -17.34519731065735
$ws.Range("C21").Value = 1.026923804045602
$ws.Range("D21").Value = -17.34519731065735
$ws.Range("E21").Value = 1.91093073642471
$ws.Range("F21").Value = -17.34519731065735
$ws.Range("G21").Value = 2.581111715860263
$ws.Range("H21").Value = 2.169333760205156
$ws.Range("I21").Value = -17.34519731065735
$ws.Range("J21").Value = -17.34519731065735
$ws.Range("K21").Value = -17.34519731065735
